$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.689.12"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'2.469.39"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'320.42"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "'92.07"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'32.85"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "'0.0852"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "'2.853.19"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "'6.86"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'15.46"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "'2.472.69"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'41.658.17"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'71.99"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'11.19"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "'239.49"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'1.93"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D27").Value = "'24.73"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").Value = "'9.69"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "'36.07"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "'155.14"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'5.41"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'0.0762"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").Value = "'16.99"
$ws.Range("E36").Value = "  -3.00%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.89"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.83"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'3.97"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'2.33"
$ws.Range("E42").Value = "  -7.17%  "
$ws.Range("D43").Value = "'1.998.45"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'18.61"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'2.94"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("D48").Value = "'2.735.76"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").Value = "'97.10"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'75.64"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").Value = "'66.80"
$ws.Range("E51").Value = "  -0.22%  "
